$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G9").Value = "ERP user/ pass: superuser/ superuser"
$ws.Range("B10").Value = "ERP user can login to ERP system"
$ws.Range("B11").Value = "Login success to ERP"

$ws.Range("B11:K11").Select()
